$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'249.41"
$cell.Style = "Normal"

$cell = $ws.Range("D4")
$cell.Value = "'5.441"
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.Value = "'0.05692"
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.Value = "'0.8078"
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.Value = "'1.032"
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.Value = "'0.1461"
$cell.Style = "Normal"

$cell = $ws.Range("D10")
$cell.Value = "'0.07806"
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.Value = "'0.03177"
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.Value = "'0.03065"
$cell.Style = "Normal"

$cell = $ws.Range("D13")
$cell.Value = "'0.09271"
$cell.Style = "Normal"

$cell = $ws.Range("D14")
$cell.Value = "'3.564"
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.Value = "'0.001646"
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.Value = "'0.04719"
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.Value = "'0.0005856"
$cell.Style = "Normal"

$cell = $ws.Range("D18")
$cell.Value = "'0.006353"
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.Value = "'0.005038"
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.Value = "'0.001042"
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.Value = "'0.0001500"
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.Value = "'0.0003203"
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.Value = "'6.425"
$cell.Style = "Normal"

$cell = $ws.Range("D25")
$cell.Value = "'2.170"
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.Value = "'0.3307"
$cell.Style = "Normal"

$cell = $ws.Range("D27")
$cell.Value = "'0.1306"
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.Value = "'0.04092"
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.Value = "'0.006915"
$cell.Style = "Normal"

$cell = $ws.Range("D43")
$cell.Value = "'0.002971"
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.Value = "'0.007772"
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.Value = "'0.00005907"
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.Value = "'0.00000000751"
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.Value = "'0.0005505"
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.Value = "'0.6831"
$cell.Style = "Normal"

$cell = $ws.Range("D49")
$cell.Value = "'0.008967"
$cell.Style = "Normal"

$cell = $ws.Range("D50")
$cell.Value = "'0.00002102"
$cell.Style = "Normal"
